$d = $word.ActiveDocument

# The <id>...</id> markers in this document were previously split across
# several runs (e.g. "<id>", "p02", "3", "r_1", "</id>"). Re-typing the
# same visible text via Find & Replace collapses each match back into a
# single run (using the formatting of the first run in the match), which
# is exactly what the commit's diff shows for the two paragraphs below.

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# <id>p023r_1</id> -> merge "<id>", "p02", "3", "r_1", "</id>" into one run
$find.Execute("<id>p023r_1</id>", $false, $false, $false, $false, $false, `
               $true, 1, $false, "<id>p023r_1</id>", 2) | Out-Null

# <id>p023v_1</id> -> merge "<id>", "p023v_1", "</id>" into one run
$find.Execute("<id>p023v_1</id>", $false, $false, $false, $false, $false, `
               $true, 1, $false, "<id>p023v_1</id>", 2) | Out-Null
